$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.28%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.05%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.112"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.28%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08003"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.22%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.412"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "28.11%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.811"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.806"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.84%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9209"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.58%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1733"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.55%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07293"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.10%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08536"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.82%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.69%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09979"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.55%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001503"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.19%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006023"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.84%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.512"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.95%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.249"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.36%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3284"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.68%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.627"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1621"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04626"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001250"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.62%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004431"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.38%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001203"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.37%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003435"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "83.35%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.74%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04468"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.16%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006963"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.36%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1343"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.05%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002245"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.65%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009825"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006615"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.56%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.23%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8206"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-57.27%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.005238"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-48.72%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002105"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.23%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002005"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.30%"
